$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 49.846522874245437
$ws.Range("C2").Value = 54.726252459006133
$ws.Range("D2").Value = 52.598329430698286
$ws.Range("E2").Value = 58.511394789979718

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 44.659973050356776
$ws.Range("C3").Value = 49.233292716779907
$ws.Range("D3").Value = 50.060005343044153
$ws.Range("E3").Value = 54.754144567851213

# Update the selection to match the new used range of interest
$ws.Range("B1:E3").Select() | Out-Null
